$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name on both sheets (remove space, add hyphen after 773)
$wsInput.Range("B1").Value = "773-RBI-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"
$wsOutput.Range("B1").Value = "773-RBI-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Update interestcalculationperiod value from "Daily" to "Same as repayment period"
$wsInput.Range("B18").Value = "Same as repayment period"

# Update selection on ProductLoanInput sheet to B18
$wsInput.Activate()
$wsInput.Range("B18").Select()

# Update selection on ProductLoanOutput sheet to B1
$wsOutput.Activate()
$wsOutput.Range("B1").Select()

# Re-activate input sheet (it is the tab that was selected originally)
$wsInput.Activate()
